$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (shared-string text changes) ---
$ws.Range("D1").Value = "Expected Name"
$ws.Range("D2").Value = "Test Fn Test Ln"
$ws.Range("B3").Value = "hemanthtestapa@unilogcorp.com"
$ws.Range("C3").Value = "'hemanth123"
$ws.Range("D3").Value = "Hemanth Sridhar"
$ws.Range("D4").Value = "General User User"

# --- Header row (A1:C1): drop the bottom border edge ---
$ws.Range("A1:C1").Borders.Item(9).LineStyle = -4142

# --- Remove the mailto hyperlink on B2 and reset its look to plain ---
$ws.Hyperlinks.Delete()
$ws.Range("B2").ClearFormats()

# --- B3 now holds an email address; give it a quiet "visited link" look ---
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B3").Font.Underline = $false
$ws.Range("B3").Borders.LineStyle = 1
$ws.Range("B3").Borders.Weight = 2

# --- D4 no longer needs wrapped text; reuse D2's plain bordered format ---
$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial(-4122)

# --- Update the active selection ---
$ws.Range("C4").Select()
